# Update countries & provincias Spain
# Applies refreshed COVID-19 country figures to the "Pais" sheet and
# swaps the Islas Malvinas / Montserrat rows (215/216), plus refreshes
# the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh timestamp (shared string reused verbatim) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 15:27"

# --- Swap Montserrat (row 215) and Islas Malvinas (row 216) ---
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("A216").Value = "Montserrat"

# --- Updated per-country figures (Casos totales, Nuevos casos, Casos
#     activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7638596
$ws.Range("C4").Value = 1684
$ws.Range("D4").Value = 4849539
$ws.Range("E4").Value = 2574428
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 214629

# Row 18 - Irak
$ws.Range("B18").Value = 382949
$ws.Range("C18").Value = 3808
$ws.Range("D18").Value = 312158
$ws.Range("E18").Value = 61327
$ws.Range("G18").Value = 65
$ws.Range("H18").Value = 9464

# Row 20 - Arabia Saudita
$ws.Range("B20").Value = 336766
$ws.Range("C20").Value = 379
$ws.Range("D20").Value = 322055
$ws.Range("E20").Value = 9813
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = 4898

# Row 26 - Alemania
$ws.Range("B26").Value = 302500
$ws.Range("C26").Value = 929
$ws.Range("E26").Value = 30995
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 9605

# Row 47 - Nepal
$ws.Range("B47").Value = 89263
$ws.Range("C47").Value = 2440
$ws.Range("D47").Value = 65202
$ws.Range("E47").Value = 23507
$ws.Range("G47").Value = 19
$ws.Range("H47").Value = 554

# Row 57 - Barein
$ws.Range("E57").Value = 4926
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 261

# Row 59 - Uzbekistan
$ws.Range("B59").Value = 58946
$ws.Range("C59").Value = 334
$ws.Range("D59").Value = 55633
$ws.Range("E59").Value = 2828
$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 485

# Row 76 - Serbia
$ws.Range("B76").Value = 33952
$ws.Range("C76").Value = 51
$ws.Range("E76").Value = 1660
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 756

# Row 78 - Dinamarca
$ws.Range("B78").Value = 30057
$ws.Range("C78").Value = 377
$ws.Range("D78").Value = 23122
$ws.Range("E78").Value = 6276
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 659

# Row 88 - Republica de Macedonia
$ws.Range("B88").Value = 18873
$ws.Range("C88").Value = 83
$ws.Range("D88").Value = 15487
$ws.Range("E88").Value = 2626
$ws.Range("G88").Value = 4
$ws.Range("H88").Value = 760

# Row 143 - Sri Lanka
$ws.Range("B143").Value = 3474
$ws.Range("C143").Value = 72
$ws.Range("D143").Value = 3259
$ws.Range("E143").Value = 202

# Row 148 - Islandia
$ws.Range("B148").Value = 2980
$ws.Range("C148").Value = 59
$ws.Range("D148").Value = 2300
$ws.Range("E148").Value = 670

# Row 176 - Burundi
$ws.Range("B176").Value = 514
$ws.Range("C176").Value = 1
$ws.Range("E176").Value = 41

# Row 179 - Islas Feroe
$ws.Range("B179").Value = 475
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 441
$ws.Range("E179").Value = 34

# Row 215 (now Islas Malvinas)
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Row 216 (now Montserrat)
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
